# Udacity PBI Project 1 Update
# Applies the "Step 2_Structure Combine and Clean the Data.xlsx" edit:
#  - Structure: fix a stray "s" typo to "Yes"
#  - Combine: add a leading "Query" column (dimension table name), shifting the
#    existing columns right, and clear the now-unused "Merge Dimension Tables"
#    data (column stays, but is blanked out) + add a blank "Foreign Keys" column
#  - Cleanup: fill in the "Address Data Gaps" column for row 4, and fully
#    populate rows 5-6 (Products / Forex) that were previously half-empty
#  - Makes "Combine" the active tab/sheet instead of "Cleanup"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Structure": fix the stray "s" -> "Yes" typo in the Review Data Types
# column for the SSBC Product Offerings row, and move the active selection.
# ---------------------------------------------------------------------------
$wsStructure = $wb.Worksheets.Item("Structure")
$wsStructure.Range("C5").Value = "Yes"
[void]$wsStructure.Range("C6").Select()

# ---------------------------------------------------------------------------
# Sheet "Combine": insert a new first column ("Query") with the dimension
# table name for each source file, pushing the old data (Merge Dimension
# Tables / Append Files / Foreign Keys) one column to the right. The "Merge
# Dimension Tables" data is cleared out (column becomes empty), and a new
# empty "Foreign Keys" column is added at the end.
# ---------------------------------------------------------------------------
$wsCombine = $wb.Worksheets.Item("Combine")

# Wipe the old A1:D6 grid contents (keeps per-cell formatting, e.g. bold
# header row, intact) before laying out the new A1:E6 grid.
$wsCombine.Range("A1:D6").ClearContents()

$wsCombine.Cells.Item(1,1).Value = "Query"
$wsCombine.Cells.Item(1,2).Value = "File"
$wsCombine.Cells.Item(1,3).Value = "Merge Dimension Tables"
$wsCombine.Cells.Item(1,4).Value = "Append Files"
$wsCombine.Cells.Item(1,5).Value = "Foreign Keys"
$wsCombine.Range("E1").Font.Bold = $true

$wsCombine.Cells.Item(2,1).Value = "Metrics"
$wsCombine.Cells.Item(2,2).Value = "CFO Metrics Tracker.xlsx"
$wsCombine.Cells.Item(2,4).Value = "No"

$wsCombine.Cells.Item(3,1).Value = "Customers"
$wsCombine.Cells.Item(3,2).Value = "Customer List (as of FY2021).txt"
$wsCombine.Cells.Item(3,4).Value = "No"

$wsCombine.Cells.Item(4,1).Value = "Sales"
$wsCombine.Cells.Item(4,2).Value = "Monthly Sales Logs"
$wsCombine.Cells.Item(4,4).Value = "Yes"

$wsCombine.Cells.Item(5,1).Value = "Products"
$wsCombine.Cells.Item(5,2).Value = "SSBC Product Offerings.pdf"
$wsCombine.Cells.Item(5,4).Value = "No"
$wsCombine.Range("A5").Interior.Color = 65535

$wsCombine.Cells.Item(6,1).Value = "Forex"
$wsCombine.Cells.Item(6,2).Value = "USD-CAD Exchange Rates.csv"
$wsCombine.Cells.Item(6,4).Value = "No"

# Column widths are re-set to track the new (auto-fit) content widths Excel
# would have computed after the edit; the host engine snaps ColumnWidth to
# 1/6-character increments, so these are the closest achievable values.
$wsCombine.Columns.Item(1).ColumnWidth = 26.5
$wsCombine.Columns.Item(3).ColumnWidth = 20.833333333333332
$wsCombine.Columns.Item(4).ColumnWidth = 25.166666666666668
$wsCombine.Columns.Item(5).ColumnWidth = 10.5

[void]$wsCombine.Range("D14").Select()

# ---------------------------------------------------------------------------
# Sheet "Cleanup": fill in remaining gaps in the "Address Data Gaps" column
# and fully populate the Products / Forex rows that were left mostly blank.
# ---------------------------------------------------------------------------
$wsCleanup = $wb.Worksheets.Item("Cleanup")

$wsCleanup.Range("F4").Value = "?"

$wsCleanup.Range("C5").Value = "Yes - 2 rows where PKProductID is null"
$wsCleanup.Range("D5").Value = "Consider drop Description?"
$wsCleanup.Range("E5").Value = "None"
$wsCleanup.Range("F5").Value = "?"

$wsCleanup.Range("C6").Value = "No"
$wsCleanup.Range("D6").Value = "Yes in import some unlabelled"
$wsCleanup.Range("E6").Value = "None"
$wsCleanup.Range("F6").Value = "?"

$wsCleanup.Columns.Item(3).ColumnWidth = 32.0
$wsCleanup.Columns.Item(4).ColumnWidth = 25.666666666666668

[void]$wsCleanup.Range("B1:B6").Select()

# ---------------------------------------------------------------------------
# Make "Combine" the active sheet/tab (it was "Cleanup" before the edit).
# ---------------------------------------------------------------------------
[void]$wsCombine.Activate()
